$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "AA2121611C2T"
$ws.Range("H2").Value = "28 jun. 2023, 14:44:44"

$ws.Range("E12").Select()
